# Updates cryptos list figures (price/volume columns) per the Jan 17 2024 GitHub Actions refresh.
# Source diff only rewrites the inline-string <t> contents of columns B-E for rows 2-51;
# row 31/32 additionally swap which coin (Filecoin/Monero) occupies each row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new value. Cells whose new text is a bare decimal number (no thousands dots) would
# otherwise be auto-coerced to a Number by Excel and lose formatting (e.g. "5.80" -> 5.8), so
# those are stamped with a Text number format first to keep them as literal strings.
$updates = [ordered]@{
    'D2' = '42.704.14'
    'E2' = '  -1.83%  '
    'D3' = '2.536.73'
    'E3' = '  -2.84%  '
    'E4' = '  -0.03%  '
    'D5' = '309.81'
    'E5' = '  -2.28%  '
    'D6' = '101.64'
    'E6' = '  +3.78%  '
    'E7' = '  -1.32%  '
    'E8' = '  +0.05%  '
    'E9' = '  -2.74%  '
    'D10' = '36.34'
    'E10' = '  +0.64%  '
    'E11' = '  -1.54%  '
    'E12' = '  -2.81%  '
    'E13' = '  +0.12%  '
    'D14' = '2.931.21'
    'E14' = '  -2.65%  '
    'D15' = '15.73'
    'D16' = '2.463.17'
    'E16' = '  -7.69%  '
    'D17' = '0.811'
    'E17' = '  -4.87%  '
    'D18' = '42.696.35'
    'E18' = '  -2.14%  '
    'D19' = '6.74'
    'E19' = '  -2.47%  '
    'D20' = '0.0₃0952'
    'E20' = '  -1.92%  '
    'E21' = '  -3.57%  '
    'D22' = '69.56'
    'E22' = '  -0.96%  '
    'D23' = '244.91'
    'E23' = '  -4.40%  '
    'E24' = '  -3.07%  '
    'E25' = '  -1.99%  '
    'E26' = '  +0.02%  '
    'D27' = '26.12'
    'E27' = '  -5.59%  '
    'D28' = '2.34'
    'E28' = '  -3.89%  '
    'D29' = '39.24'
    'E29' = '  -3.19%  '
    'E30' = '  -1.68%  '
    'B31' = 'Monero'
    'C31' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D31' = '157.05'
    'E31' = '  -0.36%  '
    'B32' = 'Filecoin'
    'C32' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D32' = '5.80'
    'E32' = '  -2.22%  '
    'E33' = '  +9.61%  '
    'E34' = '  -2.45%  '
    'E35' = '  -2.71%  '
    'D36' = '2.04'
    'E36' = '  -6.69%  '
    'D37' = '18.20'
    'E37' = '  -3.80%  '
    'E38' = '  -9.90%  '
    'E39' = '  -0.46%  '
    'E40' = '  -0.15%  '
    'D41' = '4.31'
    'E41' = '  +6.57%  '
    'D42' = '22.23'
    'E42' = '  -3.42%  '
    'E43' = '  +0.08%  '
    'E44' = '  +0.87%  '
    'D45' = '0.0301'
    'E45' = '  -1.89%  '
    'D46' = '1.985.97'
    'E46' = '  -1.91%  '
    'D47' = '8.90'
    'E47' = '  -1.60%  '
    'D48' = '2.785.46'
    'E48' = '  -2.33%  '
    'D49' = '80.83'
    'E49' = '  -3.90%  '
    'E50' = '  -1.33%  '
    'D51' = '0.851'
    'E51' = '  +7.33%  '
}

foreach ($cellRef in $updates.Keys) {
    $newValue = $updates[$cellRef]
    $range = $ws.Range($cellRef)
    if ($cellRef -match "^D" -and $newValue -match "^[+-]?\d+(\.\d+)?$") {
        $range.NumberFormat = "@"
    }
    $range.Value = $newValue
}
